# Update TPM-derived NATMI ligand-receptor statistics on the active sheet
# (Nampt-Insr.xlsx) to reflect newly computed TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.547501666666666
$ws.Range("H2").Value = 13.642505
$ws.Range("I2").Value = 0.2029775505051628
$ws.Range("J2").Value = 0.2029775505051628
$ws.Range("M2").Value = 7.655977
$ws.Range("N2").Value = 22.967931
$ws.Range("O2").Value = 0.2994795900616967
$ws.Range("P2").Value = 0.2994795900616967
$ws.Range("Q2").Value = 34.81556816746166
$ws.Range("R2").Value = 313.340113507155
$ws.Range("S2").Value = 0.0607876336170135
$ws.Range("T2").Value = 0.06078763361701349
$ws.Range("G3").Value = 4.547501666666666
$ws.Range("H3").Value = 13.642505
$ws.Range("I3").Value = 0.2029775505051628
$ws.Range("J3").Value = 0.2029775505051628
$ws.Range("O3").Value = 0.3140620915319453
$ws.Range("P3").Value = 0.3140620915319453
$ws.Range("Q3").Value = 36.51083586127999
$ws.Range("R3").Value = 328.59752275152
$ws.Range("S3").Value = 0.0637475540456825
$ws.Range("T3").Value = 0.0637475540456825
$ws.Range("G4").Value = 4.547501666666666
$ws.Range("H4").Value = 13.642505
$ws.Range("I4").Value = 0.2029775505051628
$ws.Range("J4").Value = 0.2029775505051628
$ws.Range("M4").Value = 9.879524666666667
$ws.Range("N4").Value = 29.638574
$ws.Range("O4").Value = 0.386458318406358
$ws.Range("P4").Value = 0.386458318406358
$ws.Range("Q4").Value = 44.92715488754111
$ws.Range("R4").Value = 404.34439398787
$ws.Range("S4").Value = 0.07844236284246683
$ws.Range("T4").Value = 0.07844236284246682
$ws.Range("I5").Value = 0.4559250534998339
$ws.Range("J5").Value = 0.4559250534998338
$ws.Range("M5").Value = 7.655977
$ws.Range("N5").Value = 22.967931
$ws.Range("O5").Value = 0.2994795900616967
$ws.Range("P5").Value = 0.2994795900616967
$ws.Range("Q5").Value = 78.20219398584834
$ws.Range("R5").Value = 703.8197458726351
$ws.Range("S5").Value = 0.1365402481209874
$ws.Range("T5").Value = 0.1365402481209874
$ws.Range("I6").Value = 0.4559250534998339
$ws.Range("J6").Value = 0.4559250534998338
$ws.Range("O6").Value = 0.3140620915319453
$ws.Range("P6").Value = 0.3140620915319453
$ws.Range("Q6").Value = 82.01007821776
$ws.Range("R6").Value = 738.09070395984
$ws.Range("S6").Value = 0.1431887758839719
$ws.Range("T6").Value = 0.1431887758839719
$ws.Range("I7").Value = 0.4559250534998339
$ws.Range("J7").Value = 0.4559250534998338
$ws.Range("M7").Value = 9.879524666666667
$ws.Range("N7").Value = 29.638574
$ws.Range("O7").Value = 0.386458318406358
$ws.Range("P7").Value = 0.386458318406358
$ws.Range("Q7").Value = 100.9146846275322
$ws.Range("R7").Value = 908.23216164779
$ws.Range("S7").Value = 0.1761960294948746
$ws.Range("T7").Value = 0.1761960294948746
$ws.Range("G8").Value = 7.641933666666667
$ws.Range("H8").Value = 22.925801
$ws.Range("I8").Value = 0.3410973959950033
$ws.Range("J8").Value = 0.3410973959950033
$ws.Range("M8").Value = 7.655977
$ws.Range("N8").Value = 22.967931
$ws.Range("O8").Value = 0.2994795900616967
$ws.Range("P8").Value = 0.2994795900616967
$ws.Range("Q8").Value = 58.50646838752566
$ws.Range("R8").Value = 526.558215487731
$ws.Range("S8").Value = 0.1021517083236958
$ws.Range("T8").Value = 0.1021517083236958
$ws.Range("G9").Value = 7.641933666666667
$ws.Range("H9").Value = 22.925801
$ws.Range("I9").Value = 0.3410973959950033
$ws.Range("J9").Value = 0.3410973959950033
$ws.Range("O9").Value = 0.3140620915319453
$ws.Range("P9").Value = 0.3140620915319453
$ws.Range("Q9").Value = 61.355312481056
$ws.Range("R9").Value = 552.1978123295039
$ws.Range("S9").Value = 0.1071257616022909
$ws.Range("T9").Value = 0.1071257616022909
$ws.Range("G10").Value = 7.641933666666667
$ws.Range("H10").Value = 22.925801
$ws.Range("I10").Value = 0.3410973959950033
$ws.Range("J10").Value = 0.3410973959950033
$ws.Range("M10").Value = 9.879524666666667
$ws.Range("N10").Value = 29.638574
$ws.Range("O10").Value = 0.386458318406358
$ws.Range("P10").Value = 0.386458318406358
$ws.Range("Q10").Value = 75.49867216086378
$ws.Range("R10").Value = 679.488049447774
$ws.Range("S10").Value = 0.1318199260690166
$ws.Range("T10").Value = 0.1318199260690166
Write-Output "Updated NATMI TPM values for Nampt-Insr sheet"
